# Updates the crafting-profit metric columns (H:N) on every Leve sheet
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the latest pulled
# marketboard averages, per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = -36
$ws.Range("N4").Value = -978
$ws.Range("H38").Value = 122.4
$ws.Range("I38").Value = 36
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 108
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = 264
$ws.Range("N38").Value = -3444
$ws.Range("H39").Value = 898.4706
$ws.Range("I39").Value = 109.25
$ws.Range("J39").Value = 1600
$ws.Range("K39").Value = 327.75
$ws.Range("L39").Value = 4800
$ws.Range("M39").Value = -31.75
$ws.Range("N39").Value = -5392
$ws.Range("H43").Value = 300.16666
$ws.Range("I43").Value = 250.5
$ws.Range("K43").Value = 250.5
$ws.Range("M43").Value = -181.5
$ws.Range("H58").Value = 2000.3334
$ws.Range("J58").Value = 2999.8
$ws.Range("L58").Value = 8999.400000000001
$ws.Range("N58").Value = -9299.400000000001
$ws.Range("H106").Value = 9010801
$ws.Range("I106").Value = 12346800
$ws.Range("K106").Value = 12346800
$ws.Range("M106").Value = -12346169
$ws.Range("H107").Value = 406.31818
$ws.Range("J107").Value = 119.2
$ws.Range("L107").Value = 119.2
$ws.Range("N107").Value = -3959.2
$ws.Range("H129").Value = 190418.98
$ws.Range("I129").Value = 336.2857
$ws.Range("J129").Value = 219344.61
$ws.Range("K129").Value = 1008.8571
$ws.Range("L129").Value = 658033.83
$ws.Range("M129").Value = 3991.1429
$ws.Range("N129").Value = -668033.83
$ws.Range("H132").Value = 2251.9268
$ws.Range("I132").Value = 2283.225
$ws.Range("K132").Value = 6849.674999999999
$ws.Range("M132").Value = -4319.674999999999
$ws.Range("H137").Value = 1958.303
$ws.Range("I137").Value = 1672.96
$ws.Range("K137").Value = 5018.88
$ws.Range("M137").Value = -2468.88
$ws.Range("H138").Value = 2059.2246
$ws.Range("I138").Value = 1061.7576
$ws.Range("J138").Value = 2565.6309
$ws.Range("K138").Value = 3185.2728
$ws.Range("L138").Value = 7696.8927
$ws.Range("M138").Value = 1954.7272
$ws.Range("N138").Value = -17976.8927

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4598.299
$ws.Range("I32").Value = 4421.206
$ws.Range("K32").Value = 4421.206
$ws.Range("M32").Value = -4134.206
$ws.Range("H61").Value = 1688.7037
$ws.Range("I61").Value = 1466.5
$ws.Range("J61").Value = 3466.3333
$ws.Range("K61").Value = 1466.5
$ws.Range("L61").Value = 3466.3333
$ws.Range("M61").Value = -1254.5
$ws.Range("N61").Value = -3890.3333
$ws.Range("H74").Value = 30304920
$ws.Range("I74").Value = 35716428
$ws.Range("K74").Value = 35716428
$ws.Range("M74").Value = -35715554
$ws.Range("H77").Value = 30304920
$ws.Range("I77").Value = 35716428
$ws.Range("K77").Value = 178582140
$ws.Range("M77").Value = -178577772
$ws.Range("H102").Value = 1173.2727
$ws.Range("I102").Value = 950.5714
$ws.Range("J102").Value = 1563
$ws.Range("K102").Value = 950.5714
$ws.Range("L102").Value = 1563
$ws.Range("M102").Value = 671.4286
$ws.Range("N102").Value = -4807
$ws.Range("H136").Value = 1688.7037
$ws.Range("I136").Value = 1466.5
$ws.Range("J136").Value = 3466.3333
$ws.Range("K136").Value = 4399.5
$ws.Range("L136").Value = 10398.9999
$ws.Range("M136").Value = -1849.5
$ws.Range("N136").Value = -15498.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3576.8462
$ws.Range("I105").Value = 3762.5
$ws.Range("J105").Value = 3279.8
$ws.Range("K105").Value = 3762.5
$ws.Range("L105").Value = 3279.8
$ws.Range("M105").Value = -2015.5
$ws.Range("N105").Value = -6773.8
$ws.Range("H107").Value = 1609.8
$ws.Range("I107").Value = 1475
$ws.Range("J107").Value = 1699.6666
$ws.Range("K107").Value = 1475
$ws.Range("L107").Value = 1699.6666
$ws.Range("M107").Value = 445
$ws.Range("N107").Value = -5539.6666
$ws.Range("H134").Value = 4761.485
$ws.Range("I134").Value = 4901.0356
$ws.Range("K134").Value = 14703.1068
$ws.Range("M134").Value = -12168.1068
$ws.Range("H140").Value = 42390
$ws.Range("J140").Value = 42390
$ws.Range("L140").Value = 42390
$ws.Range("N140").Value = -52750

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18275.637
$ws.Range("I31").Value = 33496.4
$ws.Range("J31").Value = 5591.6665
$ws.Range("K31").Value = 33496.4
$ws.Range("L31").Value = 5591.6665
$ws.Range("M31").Value = -33201.4
$ws.Range("N31").Value = -6181.6665
$ws.Range("H34").Value = 18275.637
$ws.Range("I34").Value = 33496.4
$ws.Range("J34").Value = 5591.6665
$ws.Range("K34").Value = 33496.4
$ws.Range("L34").Value = 5591.6665
$ws.Range("M34").Value = -33294.4
$ws.Range("N34").Value = -5995.6665
$ws.Range("H132").Value = 13291.511
$ws.Range("I132").Value = 16252.294
$ws.Range("J132").Value = 4140
$ws.Range("K132").Value = 48756.882
$ws.Range("L132").Value = 12420
$ws.Range("M132").Value = -46226.882
$ws.Range("N132").Value = -17480

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 9000
$ws.Range("J102").Value = 9000
$ws.Range("L102").Value = 27000
$ws.Range("N102").Value = -31868
$ws.Range("H113").Value = 6968.625
$ws.Range("I113").Value = 13064.875
$ws.Range("J113").Value = 872.375
$ws.Range("K113").Value = 39194.625
$ws.Range("L113").Value = 2617.125
$ws.Range("M113").Value = -37024.625
$ws.Range("N113").Value = -6957.125
$ws.Range("H131").Value = 785.03
$ws.Range("J131").Value = 806.80206
$ws.Range("L131").Value = 2420.40618
$ws.Range("N131").Value = -12500.40618

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4666666.5
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 4166666.8
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 4166666.8
$ws.Range("M7").Value = -4999888
$ws.Range("N7").Value = -4166890.8
$ws.Range("H8").Value = 4666666.5
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 4166666.8
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 4166666.8
$ws.Range("M8").Value = -4999861
$ws.Range("N8").Value = -4166944.8
$ws.Range("H122").Value = 41667664
$ws.Range("I122").Value = 16667426
$ws.Range("J122").Value = 83334740
$ws.Range("K122").Value = 50002278
$ws.Range("L122").Value = 250004220
$ws.Range("M122").Value = -49999828
$ws.Range("N122").Value = -250009120
$ws.Range("H126").Value = 3874
$ws.Range("J126").Value = 5047.4287
$ws.Range("L126").Value = 15142.2861
$ws.Range("N126").Value = -20082.2861

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 3000
$ws.Range("J14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("N14").Value = -3344
$ws.Range("H22").Value = 1782.2941
$ws.Range("I22").Value = 1500.8889
$ws.Range("J22").Value = 2098.875
$ws.Range("K22").Value = 1500.8889
$ws.Range("L22").Value = 2098.875
$ws.Range("M22").Value = -1205.8889
$ws.Range("N22").Value = -2688.875
$ws.Range("H27").Value = 1782.2941
$ws.Range("I27").Value = 1500.8889
$ws.Range("J27").Value = 2098.875
$ws.Range("K27").Value = 1500.8889
$ws.Range("L27").Value = 2098.875
$ws.Range("M27").Value = -1393.8889
$ws.Range("N27").Value = -2312.875
$ws.Range("H61").Value = 5056.067
$ws.Range("I61").Value = 2544.1
$ws.Range("K61").Value = 2544.1
$ws.Range("M61").Value = -2342.1
$ws.Range("H113").Value = 5056.067
$ws.Range("I113").Value = 2544.1
$ws.Range("K113").Value = 2544.1
$ws.Range("M113").Value = -374.0999999999999
$ws.Range("H132").Value = 1802.3636
$ws.Range("I132").Value = 1281.9565
$ws.Range("J132").Value = 2999.3
$ws.Range("K132").Value = 3845.8695
$ws.Range("L132").Value = 8997.900000000001
$ws.Range("M132").Value = -1315.8695
$ws.Range("N132").Value = -14057.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 3290
$ws.Range("J21").Value = 3290
$ws.Range("L21").Value = 3290
$ws.Range("N21").Value = -3760
$ws.Range("H24").Value = 800
$ws.Range("I24").Value = 800
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 800
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -570
$ws.Range("N24").ClearContents()
$ws.Range("H30").Value = 2519.6667
$ws.Range("I30").Value = 2000
$ws.Range("J30").Value = 2779.5
$ws.Range("K30").Value = 2000
$ws.Range("L30").Value = 2779.5
$ws.Range("M30").Value = -1893
$ws.Range("N30").Value = -2993.5
$ws.Range("H35").Value = 3290
$ws.Range("J35").Value = 3290
$ws.Range("L35").Value = 3290
$ws.Range("N35").Value = -3870
$ws.Range("H62").Value = 4175
$ws.Range("I62").Value = 3901
$ws.Range("J62").Value = 4266.3335
$ws.Range("K62").Value = 3901
$ws.Range("L62").Value = 4266.3335
$ws.Range("M62").Value = -3277
$ws.Range("N62").Value = -5514.3335
$ws.Range("H65").Value = 4175
$ws.Range("I65").Value = 3901
$ws.Range("J65").Value = 4266.3335
$ws.Range("K65").Value = 19505
$ws.Range("L65").Value = 21331.6675
$ws.Range("M65").Value = -16385
$ws.Range("N65").Value = -27571.6675
$ws.Range("H107").Value = 3497136.5
$ws.Range("I107").Value = 616.25
$ws.Range("J107").Value = 9091569
$ws.Range("K107").Value = 1848.75
$ws.Range("L107").Value = 27274707
$ws.Range("M107").Value = 71.25
$ws.Range("N107").Value = -27278547
$ws.Range("H132").Value = 1091.2
$ws.Range("I132").Value = 574.0909
$ws.Range("K132").Value = 1722.2727
$ws.Range("M132").Value = 807.7273
